$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = 11700
$ws.Range("D5").Value = 7148277.97051587
$ws.Range("D6").Value = 3517249.673455395
$ws.Range("D7").Value = 10926317.89091849
$ws.Range("D8").Value = 3640201.440161694
$ws.Range("D10").Value = 537690
$ws.Range("D11").Value = 3166283.696987225
$ws.Range("D12").Value = 4040336.403926175
$ws.Range("D13").Value = 10185758.62754164
$ws.Range("D14").Value = 2696374.240314188
$ws.Range("D15").Value = 9128358.856179321
$ws.Range("D18").Value = 54998548.8
